$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows (1253-1272) to append, each tuple: destination row, source row to copy
# formatting from (so fills/borders/number formats match exactly), and the
# A:H values for fecha, hospital, camas_ocupadas_total, camas_uci_ocupadas,
# municipio, provincia, codigo_ine, observaciones.
$newRows = @(
    @{ Dst = 1253; Src = 2;   A = 43987; B = "Hospital Universitario Miguel Servet"; C = 7;    D = 4;    E = "Zaragoza"; F = "Zaragoza"; G = 50297; H = "Fuente Aragón Hoy" }
    @{ Dst = 1254; Src = 199; A = 43987; B = "Hospital Clínico Universitario";       C = 29;   D = $null; E = "Zaragoza"; F = "Zaragoza"; G = 50297; H = "Fuente Aragón Hoy" }
    @{ Dst = 1255; Src = 2;   A = 43987; B = "Hospital Royo Villanova";              C = 2;    D = $null; E = "Zaragoza"; F = "Zaragoza"; G = 50297; H = "Fuente Aragón Hoy" }
    @{ Dst = 1256; Src = 2;   A = 43987; B = "Hospital Nuestra Señora de Gracia";    C = 1;    D = $null; E = "Zaragoza"; F = "Zaragoza"; G = 50297; H = "Fuente Aragón Hoy" }
    @{ Dst = 1257; Src = 202; A = 43987; B = "Hospital General de la Defensa";       C = 1;    D = $null; E = "Zaragoza"; F = "Zaragoza"; G = 50297; H = "Fuente Aragón Hoy" }
    @{ Dst = 1258; Src = 2;   A = 43987; B = "Hospital Obispo Polanco";              C = 6;    D = $null; E = "Teruel";   F = "Teruel";   G = 44216; H = "Fuente Aragón Hoy" }
    @{ Dst = 1259; Src = 202; A = 43987; B = "Hospital de Alcañiz";                  C = 2;    D = $null; E = "Alcañiz";  F = "Teruel";   G = 44013; H = "Fuente Aragón Hoy" }
    @{ Dst = 1260; Src = 2;   A = 43987; B = "Hospital de Barbastro";                C = 8;    D = 1;    E = "Barbastro"; F = "Huesca";   G = 22048; H = "Fuente Aragón Hoy" }
    @{ Dst = 1261; Src = 202; A = 43987; B = "Hospital San Jorge";                   C = 6;    D = 1;    E = "Huesca";   F = "Huesca";   G = 22125; H = "Fuente Aragón Hoy" }
    @{ Dst = 1262; Src = 3;   A = 43987; B = "Hospital Sagrado Corazón";             C = $null; D = $null; E = "Huesca";   F = "Huesca";   G = 22125; H = "Fuente Aragón Hoy" }
    @{ Dst = 1263; Src = 2;   A = 43987; B = "Hospital Ernest Lluch";                C = 1;    D = $null; E = "Calatayud"; F = "Zaragoza"; G = 50067; H = "Fuente Aragón Hoy" }
    @{ Dst = 1264; Src = 2;   A = 43987; B = "Hospital San José";                    C = 3;    D = $null; E = "Teruel";   F = "Teruel";   G = 44216; H = "Fuente Aragón Hoy" }
    @{ Dst = 1265; Src = 202; A = 43987; B = "Hospital Ejea – Cinco Villas";         C = $null; D = $null; E = "Ejea de los Caballeros"; F = "Zaragoza"; G = 50095; H = "Fuente Aragón Hoy" }
    @{ Dst = 1266; Src = 66;  A = 43987; B = "MAZ";                                  C = $null; D = $null; E = "Zaragoza"; F = "Zaragoza"; G = 50297; H = "Fuente Aragón Hoy" }
    @{ Dst = 1267; Src = 67;  A = 43987; B = "Hospital Viamed Montecanal";           C = $null; D = $null; E = "Zaragoza"; F = "Zaragoza"; G = 50297; H = "Fuente Aragón Hoy" }
    @{ Dst = 1268; Src = 66;  A = 43987; B = "Clínica Montpellier";                  C = $null; D = $null; E = "Zaragoza"; F = "Zaragoza"; G = 50297; H = "Fuente Aragón Hoy" }
    @{ Dst = 1269; Src = 67;  A = 43987; B = "Hospital Quirón";                      C = 1;    D = $null; E = "Zaragoza"; F = "Zaragoza"; G = 50297; H = "Fuente Aragón Hoy" }
    @{ Dst = 1270; Src = 66;  A = 43987; B = "Hospital San Juan de Dios de Zaragoza"; C = $null; D = $null; E = "Zaragoza"; F = "Zaragoza"; G = 50297; H = "Fuente Aragón Hoy" }
    @{ Dst = 1271; Src = 67;  A = 43987; B = "Clínica Viamed Santiago";              C = $null; D = $null; E = "Huesca";   F = "Huesca";   G = 22125; H = "Fuente Aragón Hoy" }
    @{ Dst = 1272; Src = 202; A = 43987; B = "Clínica El Pilar";                     C = $null; D = $null; E = "Zaragoza"; F = "Zaragoza"; G = 50297; H = "Fuente Aragón Hoy" }
)

foreach ($row in $newRows) {
    $srcRange = $ws.Range("A$($row.Src):H$($row.Src)")
    $srcRange.Copy()
    $dstRange = $ws.Range("A$($row.Dst):H$($row.Dst)")
    $dstRange.PasteSpecial(-4122)

    $ws.Cells.Item($row.Dst, 1).Value2 = $row.A
    $ws.Cells.Item($row.Dst, 2).Value = $row.B
    if ($null -eq $row.C) { $ws.Cells.Item($row.Dst, 3).Value = $null } else { $ws.Cells.Item($row.Dst, 3).Value = $row.C }
    if ($null -eq $row.D) { $ws.Cells.Item($row.Dst, 4).Value = $null } else { $ws.Cells.Item($row.Dst, 4).Value = $row.D }
    $ws.Cells.Item($row.Dst, 5).Value = $row.E
    $ws.Cells.Item($row.Dst, 6).Value = $row.F
    $ws.Cells.Item($row.Dst, 7).Value2 = $row.G
    $ws.Cells.Item($row.Dst, 8).Value = $row.H
}

$excel.CutCopyMode = 0
